$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row was inserted at row 14 ("Fruta / hortaliza, semanal"):
# existing rows 14-39 shift down to 15-40, and the new row 14 carries the
# latest week's Espárragos sample for Provincia de Linares.
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 5
$ws.Range("B14").Value = "Macroferia Regional de Talca"
$ws.Range("C14").Value = "Maule"
$ws.Range("D14").Value = 44497
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 300000000
$ws.Range("G14").Value = "Espárragos"
$ws.Range("H14").Value = "Verde"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 900
$ws.Range("L14").Value = 900
$ws.Range("M14").Value = 900
$ws.Range("N14").Value = "$/kilo"
$ws.Range("O14").Value = "Provincia de Linares"
$ws.Range("P14").Value = 900
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
